# "Running all the test cases in Chrome"
# Update the Runmode/Results columns on the "Test Cases" sheet so that the
# previously skipped notification test cases (F1 and F3) are now run ("Y")
# and the F2 test case result changes from PASS to SKIP.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C2").Value = "Y"
$ws.Range("D3").Value = "SKIP"
$ws.Range("C4").Value = "Y"

$ws.Range("C5").Select()
